# NetLiquidity / NLQ_Data.xlsx update
# - Weekly: add one new weekly data point (2023-07-12)
# - Resampled2Daily: append 11 new daily rows (2023-07-10 .. 2023-07-20)
# - Daily_TGAData: correct the last 3 existing rows (2023-07-07/08/09) and
#   append 11 new daily rows (2023-07-10 .. 2023-07-20)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Weekly": append row 94
# ---------------------------------------------------------------------
$wsWeekly = $wb.Worksheets.Item("Weekly")

$lastRowWeekly = 93
$wsWeekly.Range("A$lastRowWeekly").Copy()
$wsWeekly.Range("A94").PasteSpecial(-4122)

$wsWeekly.Cells.Item(94, 1).Value = 45119
$wsWeekly.Cells.Item(94, 2).Value = 5962.440000000001

# ---------------------------------------------------------------------
# Sheet "Resampled2Daily": append rows 649-659
# ---------------------------------------------------------------------
$wsResampled = $wb.Worksheets.Item("Resampled2Daily")

$wsResampled.Range("A648").Copy()
$wsResampled.Range("A649:A659").PasteSpecial(-4122)

$resampledRows = @(
    @{r=649; a=45117; b=6070.89},
    @{r=650; a=45118; b=6107.075},
    @{r=651; a=45119; b=5962.440000000001},
    @{r=652; a=45120; b=6015.154000000001},
    @{r=653; a=45121; b=6041.809000000001},
    @{r=654; a=45122; b=6041.809000000001},
    @{r=655; a=45123; b=6041.809000000001},
    @{r=656; a=45124; b=6054.264000000001},
    @{r=657; a=45125; b=6065.724000000001},
    @{r=658; a=45126; b=6049.782000000001},
    @{r=659; a=45127; b=6049.782000000001}
)

foreach ($row in $resampledRows) {
    $wsResampled.Cells.Item($row.r, 1).Value = $row.a
    $wsResampled.Cells.Item($row.r, 2).Value = $row.b
}

# ---------------------------------------------------------------------
# Sheet "Daily_TGAData": fix rows 646-648, then append rows 649-659
# ---------------------------------------------------------------------
$wsDaily = $wb.Worksheets.Item("Daily_TGAData")

$wsDaily.Cells.Item(646, 2).Value = 5971.019
$wsDaily.Cells.Item(647, 2).Value = 5971.019
$wsDaily.Cells.Item(648, 2).Value = 5971.019

$wsDaily.Range("A648").Copy()
$wsDaily.Range("A649:A659").PasteSpecial(-4122)

$dailyRows = @(
    @{r=649; a=45117; b=5974.554},
    @{r=650; a=45118; b=5979.732},
    @{r=651; a=45119; b=5959.359000000001},
    @{r=652; a=45120; b=6006.791000000001},
    @{r=653; a=45121; b=6034.524},
    @{r=654; a=45122; b=6034.524},
    @{r=655; a=45123; b=6034.524},
    @{r=656; a=45124; b=6037.184},
    @{r=657; a=45125; b=6018.542},
    @{r=658; a=45126; b=6002.6},
    @{r=659; a=45127; b=6002.6}
)

foreach ($row in $dailyRows) {
    $wsDaily.Cells.Item($row.r, 1).Value = $row.a
    $wsDaily.Cells.Item($row.r, 2).Value = $row.b
}
